# Removed Test Case Inter-Dependency
#
# The loan-product workbook name/identifier is tweaked (FLAT -> 1st) and the
# short name is switched from the numeric product id to a text code, on both
# the input and output sheets. Selection/active-sheet state is also reset so
# the output sheet (ProductLoanOutput) is the one left active/selected,
# instead of the input sheet sitting mid-scroll at B18.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Product "name" cell (row 1) on both sheets: FLAT variant -> 1st variant.
$ws1.Range("B1").Value = "2639-MS-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-ADD-VAR-INST-OVERDUE-FEE-1st"
$ws2.Range("B1").Value = "2639-MS-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-ADD-VAR-INST-OVERDUE-FEE-1st"

# Product "shortname" cell (row 2) on the input sheet: numeric id -> text code.
$ws1.Range("B2").Value = "263b"

# Reset selection on the input sheet back to B1 (was parked at B18), then
# leave the output sheet as the active/selected tab.
$ws1.Range("B1").Select()
$ws2.Activate()
